# Rename the "质控组" (QC group) to "北京组" (Beijing group) across both
# sheets of the dashboard workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1: column A ("组别") for every data row (2-5) changes from 质控组 to 北京组
$ws1.Range("A2").Value = "北京组"
$ws1.Range("A3").Value = "北京组"
$ws1.Range("A4").Value = "北京组"
$ws1.Range("A5").Value = "北京组"

# Sheet2: column A ("组别") for the single data row (2) changes the same way
$ws2.Range("A2").Value = "北京组"

# Update the stored view state to match: Sheet2 was the last sheet touched
# (selection left on I20), then Sheet1 becomes the active tab with A3:A5
# selected.
[void]$ws2.Activate()
[void]$ws2.Range("I20").Select()

[void]$ws1.Activate()
[void]$ws1.Range("A3:A5").Select()
